$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.462.77"
$ws.Range("E2").Value = "  -2.88%  "
$ws.Range("D3").Value = "3.731.70"
$ws.Range("E3").Value = "  -0.84%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'594.61"
$ws.Range("E5").Value = "  -3.27%  "
$ws.Range("D6").Value = "'167.68"
$ws.Range("E6").Value = "  -5.18%  "
$ws.Range("D7").Value = "3.734.13"
$ws.Range("E7").Value = "  -0.65%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").Value = "'0.526"
$ws.Range("E9").Value = "  -0.20%  "
$ws.Range("D10").Value = "'0.164"
$ws.Range("E10").Value = "  -0.81%  "
$ws.Range("D11").Value = "'6.18"
$ws.Range("E11").Value = "  -3.94%  "
$ws.Range("E12").Value = "  -4.30%  "
$ws.Range("E13").Value = "  -3.66%  "
$ws.Range("D14").Value = "'0.0000244"
$ws.Range("E14").Value = "  -3.94%  "
$ws.Range("D15").Value = "4.354.49"
$ws.Range("E15").Value = "  -0.83%  "
$ws.Range("D16").Value = "3.728.28"
$ws.Range("E16").Value = "  -0.90%  "
$ws.Range("D17").Value = "67.443.17"
$ws.Range("E17").Value = "  -2.99%  "
$ws.Range("E18").Value = "  -3.74%  "
$ws.Range("D19").Value = "'7.15"
$ws.Range("E19").Value = "  -5.01%  "
$ws.Range("D20").Value = "'17.30"
$ws.Range("E20").Value = "  +5.04%  "
$ws.Range("D21").Value = "'491.04"
$ws.Range("E21").Value = "  -3.38%  "
$ws.Range("D22").Value = "'9.22"
$ws.Range("E22").Value = "  -2.96%  "
$ws.Range("D23").Value = "'0.729"
$ws.Range("E23").Value = "  -0.32%  "
$ws.Range("D24").Value = "'85.30"
$ws.Range("E24").Value = "  -1.11%  "
$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").Value = "'0.0000144"
$ws.Range("E25").Value = "  +2.24%  "
$ws.Range("B26").Value = "Fetch.AI"
$ws.Range("C26").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D26").Value = "'2.35"
$ws.Range("E26").Value = "  -4.70%  "
$ws.Range("D27").Value = "'12.21"
$ws.Range("E27").Value = "  -4.70%  "
$ws.Range("D28").Value = "'10.05"
$ws.Range("E28").Value = "  -4.25%  "
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("E30").Value = "  -1.74%  "
$ws.Range("D31").Value = "'2.40"
$ws.Range("E31").Value = "  -4.81%  "
$ws.Range("D32").Value = "'7.79"
$ws.Range("E32").Value = "  -3.50%  "
$ws.Range("D33").Value = "'32.08"
$ws.Range("E33").Value = "  +3.75%  "
$ws.Range("E34").Value = "  -6.39%  "
$ws.Range("D35").Value = "'0.998"
$ws.Range("E35").Value = "  -0.16%  "
$ws.Range("E36").Value = "  -3.94%  "
$ws.Range("E37").Value = "  -5.19%  "
$ws.Range("E38").Value = "  -5.30%  "
$ws.Range("B39").Value = "TheGraph"
$ws.Range("C39").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D39").Value = "'0.325"
$ws.Range("E39").Value = "  -4.22%  "
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").Value = "'449.62"
$ws.Range("E40").Value = "  -4.11%  "
$ws.Range("D41").Value = "'49.16"
$ws.Range("E41").Value = "  -1.22%  "
$ws.Range("E42").Value = "  -4.43%  "
$ws.Range("D43").Value = "'2.81"
$ws.Range("E43").Value = "  -6.44%  "
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").Value = "'40.06"
$ws.Range("E46").Value = "  -9.08%  "
$ws.Range("D47").Value = "2.801.49"
$ws.Range("E47").Value = "  -4.69%  "
$ws.Range("D48").Value = "'140.40"
$ws.Range("E48").Value = "  +1.08%  "
$ws.Range("E49").Value = "  -3.71%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'24.10"
$ws.Range("E50").Value = "  +8.10%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "'25.50"
$ws.Range("E51").Value = "  -6.55%  "
